# Edit: "add creating additional reports in Word"
# Fills in previously-blank contractor/attendee/topic columns on the
# "опись к АО" sheet for rows 8, 11 and 15-24, and moves the active
# selection to H15.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("опись к АО")

# --- Row 8: "Корень" expense row gains a head-count, drops the old
#     "Препарат" value and gets Тема/Наименование подарка/Комментарии.
$ws.Range("H8").Value = 4
$ws.Range("K8").Value = ""
$ws.Range("L8").Value = "День конной авиации"
$ws.Range("M8").Value = "Китайская ваза 17 века династии Цинь"
$ws.Range("N8").Value = "4 вазы купил недорого"

# --- Row 11: "Круглый стол" row gets the standard "Представительские"
#     contractor-participant block.
$ws.Range("H11").Value = "Иванов И.А."
$ws.Range("I11").Value = "Менеджер"
$ws.Range("J11").Value = "Кафе Барикадная"
$ws.Range("K11").Value = "Альфазокс"
$ws.Range("L11").Value = "Обсуждение условий сотрудничества"

# --- Rows 15-24: same "Представительские - мероприятие" rows each get
#     the contractor + participant info filled in.
for ($r = 15; $r -le 24; $r++) {
    $ws.Range("G$r").Value = "ООО Таблетка"
    $ws.Range("H$r").Value = "Иванов И.А."
    $ws.Range("I$r").Value = "Менеджер"
    $ws.Range("J$r").Value = "Кафе Барикадная"
}

# --- Move the sheet's active selection from H22 to H15.
$ws.Activate()
$ws.Range("H15").Select()
